$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new "2020" column (Q), mirroring the formatting of the
# existing last data column (P) for every row in the table (rows 3-34).
$ws.Range("P3:P34").Copy()
$ws.Range("Q3:Q34").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Header cell for the new column.
$ws.Range("Q4").Value = 2020

# Data for the new Q column, row by row (mirrors the "-" placeholder used
# elsewhere in the sheet for missing values).
$values = @{
    5  = 51
    6  = 29
    7  = 22
    8  = 5
    9  = 3
    10 = 2
    11 = 15
    12 = 9
    13 = 5
    14 = "-"
    15 = "-"
    16 = "-"
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 7
    21 = 7
    22 = "-"
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 24
    27 = 10
    28 = 14
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $values.Keys) {
    $ws.Range("Q$row").Value = $values[$row]
}

# Selection moved to K18 in the saved file.
$ws.Range("K18").Select()
